# Edit script: insert 3 new weekly-report rows (dated 2021-11-10 / serial 44510)
# at the top of the "Zapallo italiano" data block that starts at row 172,
# pushing the existing rows 172-252 down to rows 175-255.
#
# This matches the target diff, which shows the dimension growing from
# A1:R252 to A1:R255 and every row from 175..255 carrying the same values
# that row N-3 carried before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 172, shifting everything
# from row 172 down through row 252 to rows 175 through 255.
$ws.Range("A172:R174").EntireRow.Insert() | Out-Null

# ---- New row 172 ----
$ws.Range("A172").Value = 6
$ws.Range("B172").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44510
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112032
$ws.Range("G172").Value = "Zapallo italiano"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 300
$ws.Range("K172").Value = 4000
$ws.Range("L172").Value = 5000
$ws.Range("M172").Value = 4600
$ws.Range("N172").Value = "`$/caja 50 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 92
$ws.Range("Q172").Value = 50
$ws.Range("R172").Value = "Hortaliza"

# ---- New row 173 ----
$ws.Range("A173").Value = 6
$ws.Range("B173").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C173").Value = "Metropolitana"
$ws.Range("D173").Value = 44510
$ws.Range("E173").Value = 13
$ws.Range("F173").Value = 100112032
$ws.Range("G173").Value = "Zapallo italiano"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 1400
$ws.Range("K173").Value = 5000
$ws.Range("L173").Value = 6000
$ws.Range("M173").Value = 5536
$ws.Range("N173").Value = "`$/caja 50 unidades"
$ws.Range("O173").Value = "Región de O'Higgins"
$ws.Range("P173").Value = 111
$ws.Range("Q173").Value = 50
$ws.Range("R173").Value = "Hortaliza"

# ---- New row 174 ----
$ws.Range("A174").Value = 6
$ws.Range("B174").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C174").Value = "Metropolitana"
$ws.Range("D174").Value = 44510
$ws.Range("E174").Value = 13
$ws.Range("F174").Value = 100112032
$ws.Range("G174").Value = "Zapallo italiano"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 200
$ws.Range("K174").Value = 5000
$ws.Range("L174").Value = 6000
$ws.Range("M174").Value = 5400
$ws.Range("N174").Value = "`$/caja 50 unidades"
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 108
$ws.Range("Q174").Value = 50
$ws.Range("R174").Value = "Hortaliza"
